$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently follows the
#    title heading ("Play Bonus Bunnies Slot for Free - Review 2021").
# ---------------------------------------------------------------------------
$metaPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "Meta description:*") {
        $metaPara = $para
        break
    }
}
if ($metaPara -ne $null) {
    $metaPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 2) Locate the closing paragraph (the old image-prompt paragraph) and insert
#    a new bold "Play Bonus Bunnies Slot for Free - Review 2021" paragraph
#    right before it.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)

$newParaXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Bonus Bunnies Slot for Free - Review 2021</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$lastPara.Range.InsertParagraphBefore()
$insertedPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$insertedPara.Range.InsertXML($newParaXml)

# ---------------------------------------------------------------------------
# 3) Replace the old "Create a feature image..." prompt text with the meta
#    description copy, keeping the paragraph's existing (italic) formatting.
# ---------------------------------------------------------------------------
$oldText = "Create a feature image for Bonus Bunnies that captures the fun and playful energy of the game while highlighting the bunny protagonist. The image should be in cartoon style and feature a Maya warrior with glasses who looks happy and inviting. The warrior should have a playful expression, and the image should convey a sense of fun and excitement. The background should feature green fields and fruit and vegetable symbols from the game, and there should be three rabbits in the image wreaking havoc. Overall, the image should be bright, colorful, and evoke a sense of whimsy and charm."
$newText = "Read our Bonus Bunnies slot review and play for free. Features, RTP, bonus rounds, and mobile compatibility discussed. Get ready to earn those carrots!"

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
